# Automatic update of files.
#
# 1) Column C ("Förändrad") is bumped from 45184 to 45186 (2023-09-15 -> 2023-09-17)
#    for every data row.
# 2) Every HYPERLINK() formula in columns S, T, V, W, X, Y that only has the
#    URL argument gets a second "friendly name" argument added, equal to the
#    row's "Beteckning" (column A) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Columns that may contain a HYPERLINK formula needing the friendly-name arg.
$linkCols = @(19, 20, 22, 23, 24, 25)  # S, T, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {

    # --- 1) bump the "Förändrad" date in column C ---
    $ws.Cells.Item($r, 3).Value = 45186

    # --- 2) add the friendly-name argument to HYPERLINK formulas ---
    $label = $ws.Cells.Item($r, 1).Text

    foreach ($col in $linkCols) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -match '^=HYPERLINK\("([^"]*)"\)$') {
                $url = $matches[1]
                $cell.Formula = '=HYPERLINK("' + $url + '", "' + $label + '")'
            }
        }
    }
}
